# Applies the crypto price/volume refresh captured in the Tue Jun  4 18:47:28 UTC 2024
# GitHub Actions commit: updated Price (D) / Volume(1h) (E) columns, and row 51
# switched from "Cosmos" to "FLOKI".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# cell -> new value; cells whose new value is a bare decimal number need
# NumberFormat "@" applied first so Excel stores/keeps them as text (matching
# the original inlineStr cells) instead of silently coercing to a Number.
$updates = [ordered]@{
    "D2" = @{ Value = "70.913.98"; ForceText = $false }
    "E2" = @{ Value = "  +2.20%  "; ForceText = $false }
    "D3" = @{ Value = "3.816.44"; ForceText = $false }
    "D4" = @{ Value = "0.999"; ForceText = $true }
    "E4" = @{ Value = "  +0.00%  "; ForceText = $false }
    "D5" = @{ Value = "669.27"; ForceText = $true }
    "D6" = @{ Value = "169.94"; ForceText = $true }
    "E6" = @{ Value = "  +2.45%  "; ForceText = $false }
    "D7" = @{ Value = "3.815.60"; ForceText = $false }
    "E7" = @{ Value = "  +0.94%  "; ForceText = $false }
    "E8" = @{ Value = "  -0.05%  "; ForceText = $false }
    "E9" = @{ Value = "  +1.32%  "; ForceText = $false }
    "D10" = @{ Value = "0.161"; ForceText = $true }
    "E10" = @{ Value = "  +0.67%  "; ForceText = $false }
    "E11" = @{ Value = "  +2.51%  "; ForceText = $false }
    "D12" = @{ Value = "6.99"; ForceText = $true }
    "E12" = @{ Value = "  +4.24%  "; ForceText = $false }
    "E13" = @{ Value = "  -0.45%  "; ForceText = $false }
    "D14" = @{ Value = "36.01"; ForceText = $true }
    "E14" = @{ Value = "  +0.89%  "; ForceText = $false }
    "D15" = @{ Value = "4.463.03"; ForceText = $false }
    "E15" = @{ Value = "  +1.01%  "; ForceText = $false }
    "D16" = @{ Value = "3.817.25"; ForceText = $false }
    "E16" = @{ Value = "  +0.37%  "; ForceText = $false }
    "D17" = @{ Value = "70.803.56"; ForceText = $false }
    "E17" = @{ Value = "  +2.05%  "; ForceText = $false }
    "D18" = @{ Value = "17.76"; ForceText = $true }
    "E18" = @{ Value = "  +0.23%  "; ForceText = $false }
    "D19" = @{ Value = "11.69"; ForceText = $true }
    "E19" = @{ Value = "  +21.51%  "; ForceText = $false }
    "D20" = @{ Value = "7.21"; ForceText = $true }
    "E20" = @{ Value = "  +1.16%  "; ForceText = $false }
    "D22" = @{ Value = "475.90"; ForceText = $true }
    "E22" = @{ Value = "  +1.65%  "; ForceText = $false }
    "D23" = @{ Value = "0.717"; ForceText = $true }
    "E23" = @{ Value = "  +1.68%  "; ForceText = $false }
    "D24" = @{ Value = "83.28"; ForceText = $true }
    "E24" = @{ Value = "  -0.03%  "; ForceText = $false }
    "D25" = @{ Value = "0.0000145"; ForceText = $true }
    "E25" = @{ Value = "  -2.36%  "; ForceText = $false }
    "D26" = @{ Value = "12.24"; ForceText = $true }
    "E26" = @{ Value = "  +1.73%  "; ForceText = $false }
    "E27" = @{ Value = "  +3.65%  "; ForceText = $false }
    "E28" = @{ Value = "  -1.30%  "; ForceText = $false }
    "E29" = @{ Value = "  +0.00%  "; ForceText = $false }
    "D30" = @{ Value = "3.969.19"; ForceText = $false }
    "E30" = @{ Value = "  +0.93%  "; ForceText = $false }
    "D31" = @{ Value = "2.87"; ForceText = $true }
    "E31" = @{ Value = "  +7.92%  "; ForceText = $false }
    "E32" = @{ Value = "  +2.82%  "; ForceText = $false }
    "D33" = @{ Value = "7.43"; ForceText = $true }
    "E33" = @{ Value = "  +2.10%  "; ForceText = $false }
    "D34" = @{ Value = "29.73"; ForceText = $true }
    "E34" = @{ Value = "  +3.11%  "; ForceText = $false }
    "E35" = @{ Value = "  +5.54%  "; ForceText = $false }
    "D36" = @{ Value = "9.20"; ForceText = $true }
    "E36" = @{ Value = "  +2.41%  "; ForceText = $false }
    "D37" = @{ Value = "3.773.65"; ForceText = $false }
    "E37" = @{ Value = "  +1.05%  "; ForceText = $false }
    "D38" = @{ Value = "1.00"; ForceText = $true }
    "E38" = @{ Value = "  -0.01%  "; ForceText = $false }
    "E39" = @{ Value = "  +0.65%  "; ForceText = $false }
    "D40" = @{ Value = "3.45"; ForceText = $true }
    "E40" = @{ Value = "  +1.35%  "; ForceText = $false }
    "D41" = @{ Value = "6.00"; ForceText = $true }
    "E41" = @{ Value = "  +3.17%  "; ForceText = $false }
    "D42" = @{ Value = "0.967"; ForceText = $true }
    "E42" = @{ Value = "  -0.19%  "; ForceText = $false }
    "E43" = @{ Value = "  -0.01%  "; ForceText = $false }
    "D44" = @{ Value = "2.11"; ForceText = $true }
    "E44" = @{ Value = "  +9.50%  "; ForceText = $false }
    "E45" = @{ Value = "  -0.01%  "; ForceText = $false }
    "D46" = @{ Value = "45.70"; ForceText = $true }
    "E46" = @{ Value = "  +5.44%  "; ForceText = $false }
    "D47" = @{ Value = "157.42"; ForceText = $true }
    "D48" = @{ Value = "48.05"; ForceText = $true }
    "E48" = @{ Value = "  +2.78%  "; ForceText = $false }
    "E49" = @{ Value = "  +0.90%  "; ForceText = $false }
    "D50" = @{ Value = "1.43"; ForceText = $true }
    "E50" = @{ Value = "  +4.20%  "; ForceText = $false }
    "B51" = @{ Value = "FLOKI"; ForceText = $false }
    "C51" = @{ Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"; ForceText = $false }
    "D51" = @{ Value = "0.000289"; ForceText = $true }
    "E51" = @{ Value = "  +3.55%  "; ForceText = $false }
}

foreach ($cellRef in $updates.Keys) {
    $update = $updates[$cellRef]
    $range = $ws.Range($cellRef)
    if ($update.ForceText) {
        $range.NumberFormat = "@"
    }
    $range.Value = $update.Value
}
